$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "307.66"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.81%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "38.61"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "8.68%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.096"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.07%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08113"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.40%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.962"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "4.59%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.189"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.50%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "7.943"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.97%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9294"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.72%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1435"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "11.23%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1952"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.86%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09104"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.01%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03510"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2.94%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09833"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.23%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001413"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.63%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005998"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.11%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.718"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-3.38%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3462"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.28%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1293"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-4.09%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.796"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.51%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2452"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "6.22%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04366"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.12%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.00%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.99%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001302"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "4.03%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02093"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "7.96%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05119"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.55%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007460"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.75%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01014"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.42%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.55%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002133"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.36%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009267"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-3.79%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006256"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.12%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.02%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003029"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-3.57%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.02%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.02%"
